$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 38.665747953482509
$ws.Range("C2").Value = 8.5149312994450668
$ws.Range("D2").Value = 13.820022612208845
$ws.Range("E2").Value = 3.0776763777496114

$ws.Range("B3").Value = 54.108535847815745
$ws.Range("C3").Value = 7.0653639578236493
$ws.Range("D3").Value = -10.900672193326187
$ws.Range("E3").Value = 9.7847658980804511

$ws.Range("B1:E3").Select() | Out-Null
